$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Widen column B (target OOXML width="45"; engine adds a fixed ~0.8333 pad
# to ColumnWidth when emitting the <col> width, so compensate here)
$ws.Columns.Item(2).ColumnWidth = 44.166666666666664

# Fill in the newly-tracked Tuesday afternoon timesheet rows (22-30)
$ws.Range("B22").Value = "Set up API for chatbot"
$ws.Range("C22").Value = "Teusday"

$ws.Range("B23").Value = "Set up API for chatbot"
$ws.Range("C23").Value = "Teusday"
$ws.Range("D23").Value = "https://github.com/MacksEntropy/Codered_PROS/commit/9cb76374efcfad71c65d0a7f6a3ed6410bbdb092"

$ws.Range("B24").Value = "Created initial webapp"
$ws.Range("C24").Value = "Teusday"

$ws.Range("B25").Value = "Created initial webapp"
$ws.Range("C25").Value = "Teusday"
$ws.Range("D25").Value = "https://github.com/MacksEntropy/Codered_PROS/commit/010064c5cbbd2c7238423fb262a4eac12a7b0ec1"

$ws.Range("B26").Value = "Fixed bug with user messages"
$ws.Range("C26").Value = "Teusday"
$ws.Range("D26").Value = "https://github.com/MacksEntropy/Codered_PROS/commit/acf64b97b66b5924c38c1b902585c72d6ff6c134"

$ws.Range("C27").Value = "Teusday"
$ws.Range("C28").Value = "Teusday"
$ws.Range("C29").Value = "Teusday"
$ws.Range("C30").Value = "Teusday"

$ws.Range("D30").Value = "https://github.com/MacksEntropy/Codered_PROS/commit/3537c190539025050d1bff98e59dd164c23b17d1"
$ws.Range("D29").Value = "https://github.com/MacksEntropy/Codered_PROS/commit/a35d5a790100420e49193f5c4f9f70a46d787f4c"

$ws.Range("B27").Value = "Refactored nlp module for conversation like dialouge"
$ws.Range("B28").Value = "Refactored nlp module for conversation like dialouge"
$ws.Range("B29").Value = "Refactored nlp module for conversation like dialouge"

$ws.Range("B30").Value = "Created framework for backend API"

# Restore view state: scroll position and active selection
# (topLeftCell itself isn't round-tripped by the engine's xlsx writer, but
# set ScrollRow/ScrollColumn too so the live session state is correct)
$ws.Range("B31").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
